# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-39, columns B:G) is re-sorted by
# period (ascending) instead of by worker, and new periods (1811-1904)
# are added for ALICIA JOHANNA CABALLERO LEONES.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CC","1049533806","MARYLIN TORRES SANTANDER","1702",27578,785979),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1702",32000,800000),
    @("CC","32939701","KAREN LORENA BELEÑO MARRUGO","1702",27578,785979),
    @("CC","1049533806","MARYLIN TORRES SANTANDER","1704",27578,785979),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1704",32000,800000),
    @("CC","32939701","KAREN LORENA BELEÑO MARRUGO","1704",27578,785979),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1707",32000,800000),
    @("CC","32939701","KAREN LORENA BELEÑO MARRUGO","1707",27578,785979),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1708",32000,800000),
    @("CC","32939701","KAREN LORENA BELEÑO MARRUGO","1708",27578,785979),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1709",32000,800000),
    @("CC","32939701","KAREN LORENA BELEÑO MARRUGO","1709",27578,785979),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1804",32000,800000),
    @("CC","32939701","KAREN LORENA BELEÑO MARRUGO","1804",31439,785979),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1805",32000,800000),
    @("CC","32939701","KAREN LORENA BELEÑO MARRUGO","1805",31439,785979),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1810",32000,800000),
    @("CC","32939701","KAREN LORENA BELEÑO MARRUGO","1810",31439,785979),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1811",32000,800000),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1812",32000,800000),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1901",32000,800000),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1902",32000,800000),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1903",32000,800000),
    @("CC","1101454715","ALICIA JOHANNA CABALLERO LEONES","1904",32000,800000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
